# Update the date line at the top of the document.
$d = $word.ActiveDocument
$d.Content.Find.Execute("2025-06-25 Wednesday", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "2025-06-26 Thursday", 2)

# Update the division problems/answers in the table. Addressing cells by
# (row, column) avoids ambiguity since some old values repeat in the table.
$t = $d.Tables.Item(1)

$updates = @(
    @{ Row = 1;  Col = 1; Text = "44÷7=6, 2" },
    @{ Row = 1;  Col = 2; Text = "92÷5=18, 2" },
    @{ Row = 1;  Col = 3; Text = "27÷3=9, 0" },
    @{ Row = 1;  Col = 4; Text = "83÷6=13, 5" },
    @{ Row = 1;  Col = 5; Text = "69÷8=8, 5" },

    @{ Row = 5;  Col = 1; Text = "63÷2=31, 1" },
    @{ Row = 5;  Col = 2; Text = "14÷5=2, 4" },
    @{ Row = 5;  Col = 3; Text = "86÷6=14, 2" },
    @{ Row = 5;  Col = 4; Text = "46÷6=7, 4" },
    @{ Row = 5;  Col = 5; Text = "51÷9=5, 6" },

    @{ Row = 9;  Col = 1; Text = "23÷2=11, 1" },
    @{ Row = 9;  Col = 2; Text = "58÷3=19, 1" },
    @{ Row = 9;  Col = 3; Text = "20÷7=2, 6" },
    @{ Row = 9;  Col = 4; Text = "87÷9=9, 6" },
    @{ Row = 9;  Col = 5; Text = "47÷3=15, 2" },

    @{ Row = 13; Col = 1; Text = "97÷7=13, 6" },
    @{ Row = 13; Col = 2; Text = "37÷8=4, 5" },
    @{ Row = 13; Col = 3; Text = "41÷8=5, 1" },
    @{ Row = 13; Col = 4; Text = "74÷3=24, 2" },
    @{ Row = 13; Col = 5; Text = "10÷3=3, 1" },

    @{ Row = 17; Col = 1; Text = "91÷2=45, 1" },
    @{ Row = 17; Col = 2; Text = "56÷2=28, 0" },
    @{ Row = 17; Col = 3; Text = "60÷5=12, 0" },
    @{ Row = 17; Col = 4; Text = "61÷5=12, 1" },
    @{ Row = 17; Col = 5; Text = "57÷8=7, 1" }
)

foreach ($u in $updates) {
    $cell = $t.Cell($u.Row, $u.Col)
    $cell.Range.Text = $u.Text
}
